# ------------------------------------------------------------------
# Checkpoint - Coding running well
# Update statistical results (power vs intensity re-analysis) in sheet2
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# --- Rename the data sheet to reflect the new recording id ---
$ws = $wb.Worksheets.Item("19812000")
$ws.Name = "19805000"

# --- Relabel "intensity" terminology to "power" ---
$ws.Range("E1").Value = "Power (mV^2/s), median"
$ws.Range("F1").Value = "Power (mV^2/s), IQR"
$ws.Range("A20").Value = "one-way ANOVA, power"
$ws.Range("A26").Value = "Multiple Comparison (Tukey-Kramer method), power"

# --- Updated statistic values from the re-run analysis ---
$ws.Range("B2").Value = [double]"13.349005000000034"
$ws.Range("C2").Value = [double]"5.0429759999997259"
$ws.Range("D2").Value = [double]"0.1346693543526121"
$ws.Range("E2").Value = [double]"295335.53033840557"
$ws.Range("F2").Value = [double]"154868.4657022437"
$ws.Range("G2").Value = [double]"0.54783149352144311"
$ws.Range("H2").Value = [double]"0.20543240654644268"
$ws.Range("I2").Value = [double]"12.716064453125"
$ws.Range("J2").Value = [double]"10.490753173828125"
$ws.Range("K2").Value = [double]"0.10428483665837576"
$ws.Range("L2").Value = [double]"14"
$ws.Range("B3").Value = [double]"63.321023999999852"
$ws.Range("C3").Value = [double]"87.460271999999918"
$ws.Range("E3").Value = [double]"130799.00983475964"
$ws.Range("F3").Value = [double]"116988.38552699925"
$ws.Range("H3").Value = [double]"7.1529707079602822E-3"
$ws.Range("I3").Value = [double]"16.848785400390625"
$ws.Range("J3").Value = [double]"5.72222900390625"
$ws.Range("L3").Value = [double]"3"
$ws.Range("B4").Value = [double]"21.63110400000005"
$ws.Range("C4").Value = [double]"3.5627520000005006"
$ws.Range("D4").Value = [double]"0.85073937132171151"
$ws.Range("E4").Value = [double]"69455.789809561044"
$ws.Range("F4").Value = [double]"4907.3794468311535"
$ws.Range("G4").Value = [double]"0.61070674819712234"
$ws.Range("H4").Value = [double]"7.8273477437384997E-4"
$ws.Range("I4").Value = [double]"10.80865478515625"
$ws.Range("J4").Value = [double]"2.3842620849609375"
$ws.Range("K4").Value = [double]"0.14941531559866775"
$ws.Range("L4").Value = [double]"5"
$ws.Range("B6").Value = [double]"4.0271940537589179E-3"
$ws.Range("C6").Value = [double]"1"
$ws.Range("D6").Value = [double]"-1"
$ws.Range("E6").Value = [double]"0.15369438002288383"
$ws.Range("F6").Value = [double]"0.64285714285714279"
$ws.Range("G6").Value = [double]"0.7142857142857143"
$ws.Range("I6").Value = [double]"0.61866473176178638"
$ws.Range("J6").Value = [double]"0.42857142857142855"
$ws.Range("K6").Value = [double]"-0.45238095238095238"
$ws.Range("B10").Value = [double]"14305.845705040392"
$ws.Range("D10").Value = [double]"7152.9228525201961"
$ws.Range("E10").Value = [double]"16.508662104615993"
$ws.Range("F10").Value = [double]"6.9944429386484559E-5"
$ws.Range("B11").Value = [double]"8232.3772415139028"
$ws.Range("C11").Value = [double]"19"
$ws.Range("D11").Value = [double]"433.28301271125804"
$ws.Range("B12").Value = [double]"22538.222946554295"
$ws.Range("C12").Value = [double]"21"
$ws.Range("C16").Value = [double]"-109.43190018609994"
$ws.Range("D16").Value = [double]"-75.788774857142613"
$ws.Range("E16").Value = [double]"-42.145649528185281"
$ws.Range("F16").Value = [double]"4.639545227236308E-5"
$ws.Range("C17").Value = [double]"-35.470401833081965"
$ws.Range("D17").Value = [double]"-7.920230857142764"
$ws.Range("E17").Value = [double]"19.629940118796441"
$ws.Range("F17").Value = [double]"0.74880708738790736"
$ws.Range("C18").Value = [double]"29.250001246734406"
$ws.Range("D18").Value = [double]"67.868543999999844"
$ws.Range("E18").Value = [double]"106.48708675326529"
$ws.Range("F18").Value = [double]"7.4209148172221795E-4"
$ws.Range("B22").Value = [double]"197212821663.99551"
$ws.Range("D22").Value = [double]"98606410831.997757"
$ws.Range("E22").Value = [double]"13.436592208563308"
$ws.Range("F22").Value = [double]"2.3086057016801487E-4"
$ws.Range("B23").Value = [double]"139434298274.97025"
$ws.Range("C23").Value = [double]"19"
$ws.Range("D23").Value = [double]"7338647277.6300125"
$ws.Range("B24").Value = [double]"336647119938.96576"
$ws.Range("C24").Value = [double]"21"
$ws.Range("C28").Value = [double]"-10407.900417154189"
$ws.Range("D28").Value = [double]"128050.24477535408"
$ws.Range("E28").Value = [double]"266508.38996786234"
$ws.Range("F28").Value = [double]"7.2881588095692007E-2"
$ws.Range("C29").Value = [double]"110605.65395859678"
$ws.Range("D29").Value = [double]"223988.2703406057"
$ws.Range("E29").Value = [double]"337370.88672261464"
$ws.Range("F29").Value = [double]"2.1553413097430418E-4"
$ws.Range("C30").Value = [double]"-62996.429761715175"
$ws.Range("D30").Value = [double]"95938.025565251621"
$ws.Range("E30").Value = [double]"254872.48089221842"
$ws.Range("F30").Value = [double]"0.29819585477497634"
$ws.Range("B34").Value = [double]"68.465184405113433"
$ws.Range("D34").Value = [double]"34.232592202556717"
$ws.Range("E34").Value = [double]"1.1948483361643591"
$ws.Range("F34").Value = [double]"0.32450041998049328"
$ws.Range("B35").Value = [double]"544.35297950576728"
$ws.Range("C35").Value = [double]"19"
$ws.Range("D35").Value = [double]"28.650156816093016"
$ws.Range("B36").Value = [double]"612.81816391088068"
$ws.Range("C36").Value = [double]"21"
$ws.Range("C40").Value = [double]"-13.669457229189122"
$ws.Range("D40").Value = [double]"-5.018304007393974"
$ws.Range("E40").Value = [double]"3.6328492144011744"
$ws.Range("F40").Value = [double]"0.32530873961456097"
$ws.Range("C41").Value = [double]"-6.6983585819798455"
$ws.Range("D41").Value = [double]"0.38602338518415102"
$ws.Range("E41").Value = [double]"7.4704053523481475"
$ws.Range("F41").Value = [double]"0.989497227549061"
$ws.Range("C42").Value = [double]"-4.5262283240569037"
$ws.Range("D42").Value = [double]"5.404327392578125"
$ws.Range("E42").Value = [double]"15.334883109213154"
$ws.Range("F42").Value = [double]"0.36948258883147278"

# --- Cells dropped from the re-run output (now blank) ---
$ws.Range("D3").ClearContents()
$ws.Range("G3").ClearContents()
$ws.Range("K3").ClearContents()
